$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) / Volume(1h) (E) columns with refreshed crypto quote data.
# Values are written with a leading apostrophe so Excel keeps them as literal
# text (matching the source data's inline-string / percent-as-text format)
# rather than auto-converting to numeric/percentage cell types, then the style
# is reset to Normal so no stray quote-prefix formatting is left on the cell.
$ws.Range("E2").Value = "'-0.17%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'26.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-0.27%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.664"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.78%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05950"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.54%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.610"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.72%"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'-1.90%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9102"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-4.16%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1382"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-1.59%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.04387"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'16.10%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07001"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-1.56%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03052"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-3.58%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09095"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-1.79%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001541"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.42%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006030"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-94.25%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006015"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.08%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.468"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.98%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.158"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.41%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-1.87%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-3.70%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1294"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.77%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.862"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.59%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04245"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.63%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001213"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.58%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004759"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'10.79%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'76.63%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'2.14%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.03772"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.46%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006250"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.90%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1096"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.002200"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'0.80%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01390"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'24.38%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005262"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-4.24%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.09%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.04300"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-51.37%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'9,876.84%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.09%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.09%"
$ws.Range("E50").Style = "Normal"
